$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# user_name_for_payment_review: "PERA" -> "KOAR TGR" (matches user_name_short in AK2)
$ws.Range("AC2").Value = "KOAR TGR"

# user_city_for_payment_review: "BEOGRAD" -> "Beograd" (matches user_city_short in AL2)
$ws.Range("AE2").Value = "Beograd"

# user_bban_for_payment_review: new BBAN value
$ws.Range("AF2").Value = "205-9001000243808-47"

# move selection to AL2 (updates view state/pane similarly to the target diff)
$ws.Range("AL2").Select()
